# Auto update Excel log
$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append new sensor events ---
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "15:12:17", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:12:19", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:12:25", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

$startRow = $proximity.Cells.Item($proximity.Rows.Count, 1).End(-4162).Row + 1
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $row = $startRow + $i
    $data = $proximityRows[$i]
    $rng = $proximity.Range($proximity.Cells.Item($row, 1), $proximity.Cells.Item($row, 6))
    $rng.NumberFormat = "@"
    $proximity.Cells.Item($row, 1).Value = $data[0]
    $proximity.Cells.Item($row, 2).Value = $data[1]
    $proximity.Cells.Item($row, 3).Value = $data[2]
    $proximity.Cells.Item($row, 4).Value = $data[3]
    $proximity.Cells.Item($row, 5).Value = $data[4]
    $proximity.Cells.Item($row, 6).Value = $data[5]
}

# --- Camera sheet: append new capture events ---
$camera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "15:12:19", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:12:26", "15:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = $camera.Cells.Item($camera.Rows.Count, 1).End(-4162).Row + 1
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $row = $startRow + $i
    $data = $cameraRows[$i]
    $rng = $camera.Range($camera.Cells.Item($row, 1), $camera.Cells.Item($row, 6))
    $rng.NumberFormat = "@"
    $camera.Cells.Item($row, 1).Value = $data[0]
    $camera.Cells.Item($row, 2).Value = $data[1]
    $camera.Cells.Item($row, 3).Value = $data[2]
    $camera.Cells.Item($row, 4).Value = $data[3]
    $camera.Cells.Item($row, 5).Value = $data[4]
    $camera.Cells.Item($row, 6).Value = $data[5]
}

$wb.Save()
